# Add the two new rows of hospital/department data to the bottom of the
# existing table (Hoja1), matching the style of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: AH582 | 2 | Cuarentena  (same style as existing data rows: centered
# horizontally and vertically)
$ws.Range("A5").Value = "AH582"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Cuarentena"
$ws.Range("A5:C5").HorizontalAlignment = -4108
$ws.Range("A5:C5").VerticalAlignment = -4108

# Row 6: ZX622 | 1 | Pedriatia  (A6/B6 keep the usual centered style, but C6
# only gets horizontal centering, introducing a new cell style)
$ws.Range("A6").Value = "ZX622"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Pedriatia"
$ws.Range("A6:B6").HorizontalAlignment = -4108
$ws.Range("A6:B6").VerticalAlignment = -4108
$ws.Range("C6").HorizontalAlignment = -4108

# Leave the selection on C6, mirroring where the author's cursor ended up.
$ws.Range("C6").Select()
